# Update "想去人数" (F) and "最低票价" (G) figures for the 江西-漫展信息 workbook.
# Both the "展览" sheet and the "全部类型" sheet carry the same table, so the
# same cell updates are applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G2").Value = 45

    $ws.Range("F3").Value = 304
    $ws.Range("G3").Value = 70

    $ws.Range("F4").Value = 351

    $ws.Range("F5").Value = 344

    $ws.Range("F6").Value = 1846

    $ws.Range("F7").Value = 77

    $ws.Range("F10").Value = 729

    $ws.Range("F12").Value = 354

    $ws.Range("F13").Value = 4405

    $ws.Range("F15").Value = 325

    $ws.Range("F16").Value = 1203

    $ws.Range("F17").Value = 513

    $ws.Range("F19").Value = 773

    $ws.Range("F21").Value = 405

    $ws.Range("F23").Value = 206
}
